$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")

# ---------------------------------------------------------------------------
# Row 13 must first copy C12's CURRENT (pre-edit) format before C12 itself
# is changed below, since the new "HPC" bold/italic/red-on-grey style for
# C13 is a variant of C12's current "monster2" style (bold red on grey)
# plus italic.
# ---------------------------------------------------------------------------
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C13").Font.Italic = $true
$ws.Range("C13").Value = "HPC"

# G13 gets a plain bold-red "Finished..." note (no fill, no special alignment)
$ws.Range("G13").Font.Bold = $true
$ws.Range("G13").Font.Color = 255
$ws.Range("G13").Value = "Finished on both monster2 and HPC, yet having different results"

# ---------------------------------------------------------------------------
# C12 adopts C19's CURRENT (pre-edit) format -- the green "running" style --
# because row 12 is now being (re-)run on HPC.
# ---------------------------------------------------------------------------
$ws.Range("C19").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C12").Value = "HPC"

# D12 now records that the run is happening again
$ws.Range("D12").Value = "re-running on HPC"

# G12 note updated
$ws.Range("G12").Value = "Previously on monster2"

# ---------------------------------------------------------------------------
# C19 adopts the plain grey "HPC" style (same as C11/C14/C16) since the run
# there is finished (D19 status cleared below).
# ---------------------------------------------------------------------------
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C19").Value = "HPC"
$ws.Range("D19").ClearContents()

# ---------------------------------------------------------------------------
# C11 drops its red "monster2" highlight and becomes plain "HPC"; its
# "why both HPC and monster2?" note is cleared.
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = "HPC"
$ws.Range("G11").ClearContents()

# ---------------------------------------------------------------------------
# Misc cosmetic/view updates
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 16.6
$ws.Activate()
$ws.Range("D6").Select() | Out-Null
